$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 5774.75
$ws.Cells.Item(100, 10).Value = 5099
$ws.Cells.Item(100, 12).Value = 5099
$ws.Cells.Item(100, 14).Value = -6181
$ws.Cells.Item(108, 8).Value = 99998.336
$ws.Cells.Item(108, 10).Value = 99998.336
$ws.Cells.Item(108, 12).Value = 99998.336
$ws.Cells.Item(108, 14).Value = -107678.336
$ws.Cells.Item(109, 8).Value = 99999
$ws.Cells.Item(109, 10).Value = 99999
$ws.Cells.Item(109, 12).Value = 99999
$ws.Cells.Item(109, 14).Value = -102773
$ws.Cells.Item(117, 8).Value = 94858.336
$ws.Cells.Item(117, 10).Value = 94858.336
$ws.Cells.Item(117, 12).Value = 94858.336
$ws.Cells.Item(117, 14).Value = -104036.336
$ws.Cells.Item(133, 8).Value = 76482
$ws.Cells.Item(133, 10).Value = 76482
$ws.Cells.Item(133, 12).Value = 76482
$ws.Cells.Item(133, 14).Value = -86602
$ws.Cells.Item(134, 8).Value = 99995
$ws.Cells.Item(134, 10).Value = 99995
$ws.Cells.Item(134, 12).Value = 99995
$ws.Cells.Item(134, 14).Value = -110135
$ws.Cells.Item(136, 8).Value = 99991
$ws.Cells.Item(136, 10).Value = 99991
$ws.Cells.Item(136, 12).Value = 99991
$ws.Cells.Item(136, 14).Value = -110191
$ws.Cells.Item(139, 8).Value = 98402
$ws.Cells.Item(139, 10).Value = 98402
$ws.Cells.Item(139, 12).Value = 98402
$ws.Cells.Item(139, 14).Value = -108682
$ws.Cells.Item(140, 8).Value = 80776
$ws.Cells.Item(140, 10).Value = 80776
$ws.Cells.Item(140, 12).Value = 80776
$ws.Cells.Item(140, 14).Value = -91136

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(52, 8).Value = 52080.727
$ws.Cells.Item(52, 10).Value = 52988.8
$ws.Cells.Item(52, 12).Value = 52988.8
$ws.Cells.Item(52, 14).Value = -53624.8
$ws.Cells.Item(61, 8).Value = 2023.5714
$ws.Cells.Item(61, 9).Value = 1333
$ws.Cells.Item(61, 10).Value = 3750
$ws.Cells.Item(61, 11).Value = 1333
$ws.Cells.Item(61, 12).Value = 3750
$ws.Cells.Item(61, 13).Value = -1121
$ws.Cells.Item(61, 14).Value = -4174
$ws.Cells.Item(104, 8).Value = 32130.666
$ws.Cells.Item(104, 10).Value = 32130.666
$ws.Cells.Item(104, 12).Value = 32130.666
$ws.Cells.Item(104, 14).Value = -39118.666
$ws.Cells.Item(115, 8).Value = 99983.336
$ws.Cells.Item(115, 10).Value = 99983.336
$ws.Cells.Item(115, 12).Value = 99983.336
$ws.Cells.Item(115, 14).Value = -103117.336
$ws.Cells.Item(118, 8).Value = 52997.145
$ws.Cells.Item(118, 10).Value = 52997.145
$ws.Cells.Item(118, 12).Value = 52997.145
$ws.Cells.Item(118, 14).Value = -56311.145
$ws.Cells.Item(121, 8).Value = 48993.176
$ws.Cells.Item(121, 10).Value = 48993.176
$ws.Cells.Item(121, 12).Value = 48993.176
$ws.Cells.Item(121, 14).Value = -52487.176
$ws.Cells.Item(132, 8).Value = 1932.5
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(136, 8).Value = 2023.5714
$ws.Cells.Item(136, 9).Value = 1333
$ws.Cells.Item(136, 10).Value = 3750
$ws.Cells.Item(136, 11).Value = 3999
$ws.Cells.Item(136, 12).Value = 11250
$ws.Cells.Item(136, 13).Value = -1449
$ws.Cells.Item(136, 14).Value = -16350

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(2, 8).Value = 23000
$ws.Cells.Item(2, 10).Value = 23000
$ws.Cells.Item(2, 12).Value = 23000
$ws.Cells.Item(2, 14).Value = -23226
$ws.Cells.Item(50, 8).Value = 48676
$ws.Cells.Item(50, 10).Value = 47998.332
$ws.Cells.Item(50, 12).Value = 47998.332
$ws.Cells.Item(50, 14).Value = -49146.332
$ws.Cells.Item(55, 8).Value = 34623
$ws.Cells.Item(55, 10).Value = 34623
$ws.Cells.Item(55, 12).Value = 34623
$ws.Cells.Item(55, 14).Value = -35169
$ws.Cells.Item(110, 8).Value = 84172.5
$ws.Cells.Item(110, 10).Value = 84172.5
$ws.Cells.Item(110, 12).Value = 84172.5
$ws.Cells.Item(110, 14).Value = -92352.5
$ws.Cells.Item(114, 8).Value = 90662.5
$ws.Cells.Item(114, 10).Value = 90662.5
$ws.Cells.Item(114, 12).Value = 90662.5
$ws.Cells.Item(114, 14).Value = -99340.5
$ws.Cells.Item(115, 8).Value = 72997.28999999999
$ws.Cells.Item(115, 10).Value = 74830
$ws.Cells.Item(115, 12).Value = 74830
$ws.Cells.Item(115, 14).Value = -77964
$ws.Cells.Item(122, 8).Value = 72822.14
$ws.Cells.Item(122, 10).Value = 72822.14
$ws.Cells.Item(122, 12).Value = 72822.14
$ws.Cells.Item(122, 14).Value = -82622.14
$ws.Cells.Item(132, 8).Value = 27917.836
$ws.Cells.Item(132, 10).Value = 27917.836
$ws.Cells.Item(132, 12).Value = 27917.836
$ws.Cells.Item(132, 14).Value = -38037.836
$ws.Cells.Item(135, 8).Value = 82000
$ws.Cells.Item(135, 10).Value = 82000
$ws.Cells.Item(135, 12).Value = 82000
$ws.Cells.Item(135, 14).Value = -92140
$ws.Cells.Item(138, 8).Value = 99765.336
$ws.Cells.Item(138, 10).Value = 99765.336
$ws.Cells.Item(138, 12).Value = 99765.336
$ws.Cells.Item(138, 14).Value = -110045.336
$ws.Cells.Item(140, 8).Value = 43498.617
$ws.Cells.Item(140, 10).Value = 43498.617
$ws.Cells.Item(140, 12).Value = 43498.617
$ws.Cells.Item(140, 14).Value = -53858.617

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(18, 8).Value = 25781
$ws.Cells.Item(18, 10).Value = 25781
$ws.Cells.Item(18, 12).Value = 25781
$ws.Cells.Item(18, 14).Value = -26241
$ws.Cells.Item(58, 8).Value = 1755
$ws.Cells.Item(58, 9).Value = 1481.75
$ws.Cells.Item(58, 11).Value = 1481.75
$ws.Cells.Item(58, 13).Value = -1278.75
$ws.Cells.Item(108, 8).Value = 57396.273
$ws.Cells.Item(108, 10).Value = 57396.273
$ws.Cells.Item(108, 12).Value = 57396.273
$ws.Cells.Item(108, 14).Value = -65076.273
$ws.Cells.Item(114, 8).Value = 39990.5
$ws.Cells.Item(114, 10).Value = 39990.5
$ws.Cells.Item(114, 12).Value = 39990.5
$ws.Cells.Item(114, 14).Value = -48668.5
$ws.Cells.Item(116, 8).Value = 89542
$ws.Cells.Item(116, 10).Value = 89542
$ws.Cells.Item(116, 12).Value = 89542
$ws.Cells.Item(116, 14).Value = -98720
$ws.Cells.Item(136, 8).Value = 1755
$ws.Cells.Item(136, 9).Value = 1481.75
$ws.Cells.Item(136, 11).Value = 4445.25
$ws.Cells.Item(136, 13).Value = -1895.25
$ws.Cells.Item(138, 8).Value = 94492
$ws.Cells.Item(138, 10).Value = 94492
$ws.Cells.Item(138, 12).Value = 94492
$ws.Cells.Item(138, 14).Value = -104772

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 1.7777778
$ws.Cells.Item(132, 8).Value = 7778.2856
$ws.Cells.Item(132, 9).Value = 2449
$ws.Cells.Item(132, 10).Value = 8666.5
$ws.Cells.Item(132, 11).Value = 22041
$ws.Cells.Item(132, 12).Value = 77998.5
$ws.Cells.Item(132, 13).Value = -19511
$ws.Cells.Item(132, 14).Value = -83058.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(108, 8).Value = 51965.637
$ws.Cells.Item(108, 10).Value = 51965.637
$ws.Cells.Item(108, 12).Value = 51965.637
$ws.Cells.Item(108, 14).Value = -59645.637
$ws.Cells.Item(109, 8).Value = 26214.666
$ws.Cells.Item(109, 10).Value = 26214.666
$ws.Cells.Item(109, 12).Value = 26214.666
$ws.Cells.Item(109, 14).Value = -28294.666
$ws.Cells.Item(114, 8).Value = 63645.184
$ws.Cells.Item(114, 10).Value = 63645.184
$ws.Cells.Item(114, 12).Value = 63645.184
$ws.Cells.Item(114, 14).Value = -72323.18400000001
$ws.Cells.Item(132, 8).Value = 3715.5
$ws.Cells.Item(132, 9).Value = 3141.6
$ws.Cells.Item(132, 10).Value = 4672
$ws.Cells.Item(132, 11).Value = 9424.799999999999
$ws.Cells.Item(132, 12).Value = 14016
$ws.Cells.Item(132, 13).Value = -6894.799999999999
$ws.Cells.Item(132, 14).Value = -19076
$ws.Cells.Item(135, 8).Value = 52291
$ws.Cells.Item(135, 10).Value = 52291
$ws.Cells.Item(135, 12).Value = 52291
$ws.Cells.Item(135, 14).Value = -62431
$ws.Cells.Item(140, 8).Value = 90430.664
$ws.Cells.Item(140, 10).Value = 90396
$ws.Cells.Item(140, 12).Value = 90396
$ws.Cells.Item(140, 14).Value = -100756

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 4499.5884
$ws.Cells.Item(16, 9).Value = 4082.3076
$ws.Cells.Item(16, 10).Value = 5855.75
$ws.Cells.Item(16, 11).Value = 4082.3076
$ws.Cells.Item(16, 12).Value = 5855.75
$ws.Cells.Item(16, 13).Value = -3912.3076
$ws.Cells.Item(16, 14).Value = -6195.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2042.7142
$ws.Cells.Item(81, 9).Value = 1883.1666
$ws.Cells.Item(81, 10).Value = 3000
$ws.Cells.Item(81, 11).Value = 3766.3332
$ws.Cells.Item(81, 12).Value = 6000
$ws.Cells.Item(81, 13).Value = -2705.3332
$ws.Cells.Item(81, 14).Value = -8122
$ws.Cells.Item(84, 8).Value = 2042.7142
$ws.Cells.Item(84, 9).Value = 1883.1666
$ws.Cells.Item(84, 10).Value = 3000
$ws.Cells.Item(84, 11).Value = 18831.666
$ws.Cells.Item(84, 12).Value = 30000
$ws.Cells.Item(84, 13).Value = -13527.666
$ws.Cells.Item(84, 14).Value = -40608
$ws.Cells.Item(96, 8).Value = 5268024.5
$ws.Cells.Item(96, 9).Value = 1239.1428
$ws.Cells.Item(96, 10).Value = 17557192
$ws.Cells.Item(96, 11).Value = 1239.1428
$ws.Cells.Item(96, 12).Value = 17557192
$ws.Cells.Item(96, 13).Value = 133.8571999999999
$ws.Cells.Item(96, 14).Value = -17559938
$ws.Cells.Item(121, 8).Value = 37354.855
$ws.Cells.Item(121, 10).Value = 37354.855
$ws.Cells.Item(121, 12).Value = 37354.855
$ws.Cells.Item(121, 14).Value = -40848.855

# --- Clear N132 on ARM (cell removed from sheet in target state) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 14).ClearContents()
